# Refresh the fantasy roster table (Oyuncu Adı / Pozisyon / Takım)
# with the new player / position / team data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A2').Value = 'Scoot Henderson'
$ws.Range('B2').Value = 'PG'
$ws.Range('C2').Value = 'Portland Trail Blazers'
$ws.Range('A3').Value = 'Isaiah Joe'
$ws.Range('B3').Value = 'PG,SG'
$ws.Range('C3').Value = 'Oklahoma City Thunder'
$ws.Range('A4').Value = 'Dalano Banton'
$ws.Range('B4').Value = 'SG,SF'
$ws.Range('C4').Value = 'Portland Trail Blazers'
$ws.Range('A5').Value = 'Alex Caruso'
$ws.Range('B5').Value = 'SG,SF'
$ws.Range('C5').Value = 'Oklahoma City Thunder'
$ws.Range('A6').Value = 'Jalen Williams'
$ws.Range('B6').Value = 'SG,SF,PF,C'
$ws.Range('C6').Value = 'Oklahoma City Thunder'
$ws.Range('A7').Value = 'Kawhi Leonard'
$ws.Range('B7').Value = 'SG,SF,PF'
$ws.Range('C7').Value = 'LA Clippers'
$ws.Range('A8').Value = 'Norman Powell'
$ws.Range('B8').Value = 'SG,SF'
$ws.Range('C8').Value = 'LA Clippers'
$ws.Range('A9').Value = 'D''Angelo Russell'
$ws.Range('B9').Value = 'PG'
$ws.Range('C9').Value = 'Brooklyn Nets'
$ws.Range('A10').Value = 'Desmond Bane'
$ws.Range('B10').Value = 'SG,SF'
$ws.Range('C10').Value = 'Memphis Grizzlies'
$ws.Range('A11').Value = 'Myles Turner'
$ws.Range('B11').Value = 'C'
$ws.Range('C11').Value = 'Indiana Pacers'
$ws.Range('A12').Value = 'Trae Young'
$ws.Range('B12').Value = 'PG'
$ws.Range('C12').Value = 'Atlanta Hawks'
$ws.Range('A13').Value = 'LeBron James'
$ws.Range('B13').Value = 'SF,PF'
$ws.Range('C13').Value = 'Los Angeles Lakers'
$ws.Range('A14').Value = 'Devin Booker'
$ws.Range('B14').Value = 'PG,SG'
$ws.Range('C14').Value = 'Phoenix Suns'
$ws.Range('A15').Value = 'Jalen Brunson'
$ws.Range('B15').Value = 'PG'
$ws.Range('C15').Value = 'New York Knicks'
$ws.Range('A16').Value = 'Walker Kessler'
$ws.Range('B16').Value = 'C'
$ws.Range('C16').Value = 'Utah Jazz'
$ws.Range('A17').Value = 'Immanuel Quickley'
$ws.Range('B17').Value = 'PG,SG'
$ws.Range('C17').Value = 'Toronto Raptors'
$ws.Range('A18').Value = 'Brandon Ingram'
$ws.Range('B18').Value = 'SG,SF,PF'
$ws.Range('C18').Value = 'New Orleans Pelicans'
$ws.Range('A19').Value = 'Jimmy Butler'
$ws.Range('B19').Value = 'SF,PF'
$ws.Range('C19').Value = 'Miami Heat'
